# 23 dec 2023 update
# Record payments received for Dec 11-17, 2023 on the
# "SM5000.1-SEPT (2)" worksheet (the workbook's active sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SM5000.1-SEPT (2)")

# Dec 11 & Dec 12, 2023 -> rows 31 & 32 (columns B:D)
$ws.Range("B31").Value = 45271
$ws.Range("C31").Value = 100
$ws.Range("D31").Value = 1

$ws.Range("B32").Value = 45272
$ws.Range("C32").Value = 100
$ws.Range("D32").Value = 1

# Dec 13-17, 2023 -> rows 3-7 (columns F:H)
$ws.Range("F3").Value = 45273
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = 1

$ws.Range("F4").Value = 45274
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = 1

$ws.Range("F5").Value = 45275
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 1

$ws.Range("F6").Value = 45276
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = 1

$ws.Range("F7").Value = 45277
$ws.Range("G7").Value = 100
$ws.Range("H7").Value = 1

# Update the saved view: scroll/selection moves to the newly entered block.
$ws.Activate()
$ws.Range("H3:H7").Select()
